$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Add term accession number + source ref examples for the Tags rows
$ws2.Range("E13").Value = "http://purl.obolibrary.org/obo/NCIT_C14258"
$ws2.Range("E14").Value = "NCIT"

$ws2.Range("B13").Copy()
$ws2.Range("E13").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("E13").WrapText = $true
$ws2.Rows.Item(13).RowHeight = 72

# Insert "Plant" term as new E12 value, shifting old E:G tags to F:H
$ws2.Range("E12").Value = "Plant"

$ws2.Activate()
$ws2.Range("H13").Select()

Write-Host "done"
